$d = $word.ActiveDocument

function Set-ParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.Text = $newText
}

Set-ParaText 4 "Unidad Económica: Financiera Independencai"
Set-ParaText 5 "Periodo: 3 mese"
Set-ParaText 6 "Fecha: 2023-02-05"

Set-ParaText 8 "legalmente_constituida: 0"
Set-ParaText 9 "convenio_cooperacion: 0"
Set-ParaText 10 "convenio_aprendizaje: 0"
Set-ParaText 11 "convenio_marco: 0"

Set-ParaText 13 "personal_capacitado: 1"
Set-ParaText 14 "areas_especializadas: 1"
Set-ParaText 15 "mentor_licenciatura: 1"
Set-ParaText 16 "plan_formacion: 1"
Set-ParaText 17 "capacidad_plan: 1"
Set-ParaText 18 "puestos_aprendizaje: 1"
